$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: PRM -> PRB, Peralatan Mandi -> Peralatan Bersih-bersih
$ws.Range("A2").Value = "PRB"
$ws.Range("B2").Value = "Peralatan Bersih-bersih"

# Delete row 3 (PMK / Peralatan Makan) entirely
$ws.Rows("3:3").Delete()

# Adjust column B width
$ws.Columns("B:B").ColumnWidth = 20

# Update selection
$ws.Range("I3").Select()
